$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ("Known issue" #2) content updates -----------------------------
# Category
$ws.Range("E4").Value = "Imperas Model"

# Example failing test(s)
$ws.Range("H4").Value = "corev_rand_pulp_instr_test"

# Example error message(s)
$errMsg = "# Info (IDV) Instruction executed prior to mismatch '0x80(_start+0): 0800006f j       100'`n" + `
          "# Error (IDV) PC mismatch (HartId:0, PC:0x00000100 _start_main+0):`n" + `
          "# Error (IDV) Mismatch 0>`n" + `
          "# Error (IDV)   . dut:0x1a110800 debug_rom+0`n" + `
          "# Error (IDV)   . ref:0x00000100 _start_main+0`n" + `
          "# Error (IDV) Insn. bit pattern mismatch (HartId:0, PC:0x00000100 _start_main+0):"
$ws.Range("I4").Value = $errMsg

# Known issue description
$ws.Range("D4").Value = "Mismatch between RTL and Imperas model on first debug request"

# Row 4 now needs much more vertical space for the wrapped error message.
$ws.Rows.Item(4).RowHeight = 115.2

# --- Selection -------------------------------------------------------------
[void]$ws.Range("H3").Select()
